# Update 2022 (column I) violent crime figures to reflect data through 2022-09-23
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5297
$ws.Range("I3").Value = 5559
$ws.Range("D4").Value = 1933
$ws.Range("H4").Value = 1674
$ws.Range("I4").Value = 1269
$ws.Range("I5").Value = 514
$ws.Range("I6").Value = 6058
$ws.Range("D7").Value = 28123
$ws.Range("H7").Value = 25984
$ws.Range("I7").Value = 18697

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I8").Value = 1121
$ws.Range("I9").Value = 88
$ws.Range("I10").Value = 135
$ws.Range("I11").Value = 278
$ws.Range("I13").Value = 34
$ws.Range("I19").Value = 512
$ws.Range("I20").Value = 448
$ws.Range("I29").Value = 1199
$ws.Range("I33").Value = 852
$ws.Range("I37").Value = 603
$ws.Range("I42").Value = 633
$ws.Range("I43").Value = 151
$ws.Range("I48").Value = 251
$ws.Range("I51").Value = 214
$ws.Range("I52").Value = 413
$ws.Range("I53").Value = 196
$ws.Range("H55").Value = 310
$ws.Range("I55").Value = 211
$ws.Range("D63").Value = 323
$ws.Range("I63").Value = 73
$ws.Range("I65").Value = 428
$ws.Range("I67").Value = 746
$ws.Range("I73").Value = 167
$ws.Range("I75").Value = 57
$ws.Range("I76").Value = 275
$ws.Range("I79").Value = 528
$ws.Range("I83").Value = 393
$ws.Range("I84").Value = 159
$ws.Range("I85").Value = 841
$ws.Range("I86").Value = 114
$ws.Range("I88").Value = 174
$ws.Range("I90").Value = 229
$ws.Range("I95").Value = 302
$ws.Range("I96").Value = 197
$ws.Range("I99").Value = 347
$ws.Range("D101").Value = 28123
$ws.Range("H101").Value = 25984
$ws.Range("I101").Value = 18697

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 233
$ws.Range("I3").Value = 335
$ws.Range("I7").Value = 841

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 151
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 413

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 278

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 317
$ws.Range("I5").Value = 33
$ws.Range("I6").Value = 361
$ws.Range("I7").Value = 1121

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 42
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 196

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 58
$ws.Range("I7").Value = 197

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 200
$ws.Range("I7").Value = 603

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I5").Value = 12
$ws.Range("I7").Value = 347

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 274
$ws.Range("I4").Value = 47
$ws.Range("I6").Value = 235
$ws.Range("I7").Value = 746

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 143
$ws.Range("I3").Value = 126
$ws.Range("I7").Value = 428

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 393

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 302

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 321
$ws.Range("I6").Value = 265
$ws.Range("I7").Value = 852

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I3").Value = 13
$ws.Range("I6").Value = 79

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 353
$ws.Range("I3").Value = 415
$ws.Range("I5").Value = 43
$ws.Range("I6").Value = 323
$ws.Range("I7").Value = 1199

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 155
$ws.Range("I7").Value = 512

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I4").Value = 26
$ws.Range("I7").Value = 251

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 125
$ws.Range("I7").Value = 275

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 218
$ws.Range("I4").Value = 47
$ws.Range("I5").Value = 23
$ws.Range("I7").Value = 633

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("I4").Value = 11
$ws.Range("I6").Value = 34

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 46
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 68
$ws.Range("H4").Value = 22
$ws.Range("H7").Value = 310
$ws.Range("I7").Value = 211

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 171
$ws.Range("I7").Value = 528

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 129
$ws.Range("I3").Value = 137
$ws.Range("I6").Value = 139
$ws.Range("I7").Value = 448

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I6").Value = 24
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 167

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 58
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 80
$ws.Range("I3").Value = 53
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 56
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 214

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I3").Value = 28
$ws.Range("I7").Value = 151
